$d = $word.ActiveDocument

# The intro paragraph used to end with a sentence ("The objective of this
# document is to enable future collaboration when new developers join the
# project in the future.") spread across four runs. That whole trailing
# thought is being dropped, leaving the paragraph ending right after
# "...things like code style. ".

# Locate the full legacy passage that must disappear.
$full = $d.Content.Duplicate
$full.Find.Text = "The objective of this document is to enable future collaboration when new developers join the project in the future."
$full.Find.Execute() | Out-Null

# Locate where the separately-authored trailing runs ("enable future
# collaboration " / "when " / "new developers...future.") begin; that
# sub-range covers three whole <w:r> elements which we can just delete.
$tail = $d.Content.Duplicate
$tail.Find.Text = "enable future collaboration when new developers join the project in the future."
$tail.Find.Execute() | Out-Null

# Delete the three trailing runs outright (they disappear completely).
$d.Range($tail.Start, $tail.End).Delete()

# What's left to trim is "The objective of this document is to " sitting at
# the tail end of the run that starts with "It also includes...code style. ".
# Nudging a formatting property on that sub-range before collapsing its text
# keeps it from being silently re-absorbed into the preceding run, so the
# surviving run boundaries stay exactly where the original authoring left
# them.
$lead = $d.Range($full.Start, $tail.Start)
$lead.Font.Bold = $true
$lead.Font.Bold = $false
$lead.Text = ""
